# Auto-generated Excel COM-interop script to apply the market-data update
# described by the Coeurl_Profits.xlsx diff (chore: update Sheets via scheduled runner).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 39663.426
$ws.Range("I74").Value = 61355.277
$ws.Range("J74").Value = 13633.2
$ws.Range("K74").Value = 61355.277
$ws.Range("L74").Value = 13633.2
$ws.Range("M74").Value = -60419.277
$ws.Range("N74").Value = -15505.2
$ws.Range("H77").Value = 39663.426
$ws.Range("I77").Value = 61355.277
$ws.Range("J77").Value = 13633.2
$ws.Range("K77").Value = 306776.385
$ws.Range("L77").Value = 68166
$ws.Range("M77").Value = -302096.385
$ws.Range("N77").Value = -77526
$ws.Range("H100").Value = 2651.5
$ws.Range("I100").Value = 2651.5
$ws.Range("K100").Value = 2651.5
$ws.Range("M100").Value = -2110.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3101.5
$ws.Range("I61").Value = 2912.3215
$ws.Range("K61").Value = 2912.3215
$ws.Range("M61").Value = -2700.3215
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H136").Value = 3101.5
$ws.Range("I136").Value = 2912.3215
$ws.Range("K136").Value = 8736.9645
$ws.Range("M136").Value = -6186.9645
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3604.2
$ws.Range("I86").Value = 2178.4783
$ws.Range("J86").Value = 20000
$ws.Range("K86").Value = 2178.4783
$ws.Range("L86").Value = 20000
$ws.Range("M86").Value = -1055.4783
$ws.Range("N86").Value = -22246
$ws.Range("H89").Value = 3604.2
$ws.Range("I89").Value = 2178.4783
$ws.Range("J89").Value = 20000
$ws.Range("K89").Value = 10892.3915
$ws.Range("L89").Value = 100000
$ws.Range("M89").Value = -5276.391500000002
$ws.Range("N89").Value = -111232
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9430.727999999999
$ws.Range("I16").Value = 340.5
$ws.Range("K16").Value = 340.5
$ws.Range("M16").Value = -53.5
$ws.Range("H22").Value = 272.27274
$ws.Range("I22").Value = 255.625
$ws.Range("K22").Value = 255.625
$ws.Range("M22").Value = 94.375
$ws.Range("H28").Value = 31520
$ws.Range("J28").Value = 31520
$ws.Range("L28").Value = 31520
$ws.Range("N28").Value = -32010
$ws.Range("H58").Value = 2074.7144
$ws.Range("I58").Value = 1789.6
$ws.Range("J58").Value = 2787.5
$ws.Range("K58").Value = 1789.6
$ws.Range("L58").Value = 2787.5
$ws.Range("M58").Value = -1586.6
$ws.Range("N58").Value = -3193.5
$ws.Range("H62").Value = 6786.857
$ws.Range("I62").Value = 9332.333000000001
$ws.Range("J62").Value = 4877.75
$ws.Range("K62").Value = 9332.333000000001
$ws.Range("L62").Value = 4877.75
$ws.Range("M62").Value = -8708.333000000001
$ws.Range("N62").Value = -6125.75
$ws.Range("H65").Value = 6786.857
$ws.Range("I65").Value = 9332.333000000001
$ws.Range("J65").Value = 4877.75
$ws.Range("K65").Value = 46661.665
$ws.Range("L65").Value = 24388.75
$ws.Range("M65").Value = -43541.665
$ws.Range("N65").Value = -30628.75
$ws.Range("H105").Value = 976.25
$ws.Range("I105").Value = 799.8333
$ws.Range("K105").Value = 799.8333
$ws.Range("M105").Value = 947.1667
$ws.Range("H113").Value = 9430.727999999999
$ws.Range("I113").Value = 340.5
$ws.Range("K113").Value = 340.5
$ws.Range("M113").Value = 1829.5
$ws.Range("H132").Value = 2846.6826
$ws.Range("J132").Value = 3519.9
$ws.Range("L132").Value = 10559.7
$ws.Range("N132").Value = -15619.7
$ws.Range("H136").Value = 2074.7144
$ws.Range("I136").Value = 1789.6
$ws.Range("J136").Value = 2787.5
$ws.Range("K136").Value = 5368.799999999999
$ws.Range("L136").Value = 8362.5
$ws.Range("M136").Value = -2818.799999999999
$ws.Range("N136").Value = -13462.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 8000
$ws.Range("I42").Value = 6000
$ws.Range("J42").Value = 10000
$ws.Range("K42").Value = 18000
$ws.Range("L42").Value = 30000
$ws.Range("M42").Value = -17466
$ws.Range("N42").Value = -31068
$ws.Range("H107").Value = 1071.0588
$ws.Range("I107").Value = 1169.3334
$ws.Range("J107").Value = 1050
$ws.Range("K107").Value = 3508.0002
$ws.Range("L107").Value = 3150
$ws.Range("M107").Value = -1588.0002
$ws.Range("N107").Value = -6990
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3481
$ws.Range("J80").Value = 4250
$ws.Range("L80").Value = 4250
$ws.Range("N80").Value = -6246
$ws.Range("H83").Value = 3481
$ws.Range("J83").Value = 4250
$ws.Range("L83").Value = 21250
$ws.Range("N83").Value = -31234
$ws.Range("H113").Value = 11002.091
$ws.Range("I113").Value = 9144.286
$ws.Range("J113").Value = 14253.25
$ws.Range("K113").Value = 9144.286
$ws.Range("L113").Value = 14253.25
$ws.Range("M113").Value = -6974.286
$ws.Range("N113").Value = -18593.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1299.8
$ws.Range("I22").Value = 999
$ws.Range("K22").Value = 999
$ws.Range("M22").Value = -704
$ws.Range("H27").Value = 1299.8
$ws.Range("I27").Value = 999
$ws.Range("K27").Value = 999
$ws.Range("M27").Value = -892
$ws.Range("H68").Value = 28889.8
$ws.Range("J68").Value = 45333
$ws.Range("L68").Value = 45333
$ws.Range("N68").Value = -46831
$ws.Range("H71").Value = 28889.8
$ws.Range("J71").Value = 45333
$ws.Range("L71").Value = 226665
$ws.Range("N71").Value = -234153
$ws.Range("H132").Value = 3603.889
$ws.Range("I132").Value = 2619.3572
$ws.Range("J132").Value = 7049.75
$ws.Range("K132").Value = 7858.071599999999
$ws.Range("L132").Value = 21149.25
$ws.Range("M132").Value = -5328.071599999999
$ws.Range("N132").Value = -26209.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16843.334
$ws.Range("I45").Value = 12000
$ws.Range("K45").Value = 12000
$ws.Range("M45").Value = -11509
$ws.Range("H54").Value = 30048.75
$ws.Range("J54").Value = 30597.5
$ws.Range("L54").Value = 30597.5
$ws.Range("N54").Value = -31637.5
$ws.Range("H100").Value = 1266.9445
$ws.Range("I100").Value = 1115
$ws.Range("J100").Value = 1505.7142
$ws.Range("K100").Value = 2230
$ws.Range("L100").Value = 3011.4284
$ws.Range("M100").Value = -1689
$ws.Range("N100").Value = -4093.4284
$ws.Range("H107").Value = 1842.5834
$ws.Range("I107").Value = 1785.5
$ws.Range("K107").Value = 5356.5
$ws.Range("M107").Value = -3436.5
$ws.Range("H136").Value = 2366.7932
$ws.Range("I136").Value = 2317.72
$ws.Range("K136").Value = 6953.16
$ws.Range("M136").Value = -4403.16
